$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats -- copy only cell formatting (styles), not values.
$xlPasteFormats = -4122

# --- 1. Duplicate the 4-row "collect card" block (rows 2-5) down through
#        row 31, renumbering column A (the "no" column) sequentially --------

# Card payload for the recurring 4-row cycle (same text as rows 2-5, i.e.
# shared-string indices 5,6,7,8 / 9,10,7,11 / 12,13,14,15 / 16,17,14,18).
$cardB = @(
    "2022/07/26 shu yamino organized neatly",
    "2022/01/23 【📻📶】first horror game stream (The Radio Station)",
    "2021/12/22 how i got here",
    "2022/04/27 how many bananas did we get (Just Chatting)"
)
$cardC = @(
    "https://www.youtube.com/watch?v=GvBPFY3altg&t=10821s",
    "https://www.youtube.com/watch?v=4TLRTBpsjmg&t=2988s",
    "https://www.youtube.com/watch?v=NkbG17qCTSA&t=7s",
    "https://www.youtube.com/watch?v=BTNvN_nK9eI&t=2553s"
)
$cardD = @("(2)", "(2)", "(3)", "(3)")
$cardE = @("spvoice_1", "spvoice_2", "spvoice_3", "spvoice_4")

# Write the values FIRST (so column A lands as a genuine number), then copy
# the template's formatting on top -- applying the "Text" number format
# (style s="1") before the value is entered would otherwise turn the number
# into a literal text string, unlike the source rows.
for ($r = 6; $r -le 31; $r++) {
    $idx = ($r - 6) % 4
    $ws.Range("A$r").Value = $r - 1
    $ws.Range("B$r").Value = $cardB[$idx]
    $ws.Range("C$r").Value = $cardC[$idx]
    $ws.Range("D$r").Value = $cardD[$idx]
    $ws.Range("E$r").Value = $cardE[$idx]
}

# Re-apply the template formatting (styles s="1"/s="2") row-block by
# row-block, always copying exactly as many source rows as the destination
# needs -- PasteSpecial always stamps the *whole* clipboard, so a
# mismatched source/destination size would overflow into later rows.
$ws.Range("A2:E5").Copy()
$ws.Range("A6:E9").PasteSpecial($xlPasteFormats)
$ws.Range("A10:E13").PasteSpecial($xlPasteFormats)
$ws.Range("A14:E17").PasteSpecial($xlPasteFormats)
$ws.Range("A18:E21").PasteSpecial($xlPasteFormats)
$ws.Range("A22:E25").PasteSpecial($xlPasteFormats)
$ws.Range("A26:E29").PasteSpecial($xlPasteFormats)

$ws.Range("A2:E3").Copy()
$ws.Range("A30:E31").PasteSpecial($xlPasteFormats)

# --- 2. Append six more blank (but styled) collect-card rows, 32-37 -------
$ws.Range("B2:E2").Copy()
$ws.Range("B32").PasteSpecial($xlPasteFormats)
$ws.Range("B33").PasteSpecial($xlPasteFormats)
$ws.Range("B34").PasteSpecial($xlPasteFormats)
$ws.Range("B35").PasteSpecial($xlPasteFormats)
$ws.Range("B36").PasteSpecial($xlPasteFormats)
$ws.Range("B37").PasteSpecial($xlPasteFormats)
